$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Fix the "CasesTab" query (cell B2): remove the erroneous trailing
# `Cohort` column from the RETURN clause (co.cohort_description), which
# was causing a query error. This is the main content fix described by
# the commit message ("Fixed variables and query errors").
# -----------------------------------------------------------------------
$caseQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Beagle']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $caseQuery

# -----------------------------------------------------------------------
# Adjust row heights (matches the re-wrapped text heights after the
# query text/zoom changes).
# -----------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# -----------------------------------------------------------------------
# Update window zoom and selection to match the saved view state.
# -----------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 130
$ws.Range("B2").Select() | Out-Null
